$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.695.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.484.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.21%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.70%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.478.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.655.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.130.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000256"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.486.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("E21").Value = "  +7.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.506"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "505.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.44%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "92.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.667.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +8.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "556.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.945"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.92%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0415"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
